$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new data row at row 9 (pushing LOLAWEST / ditol / seringat / the
#    subtotal row / the footer row all down by one), using row 8 as the
#    formatting template so fonts / fills / borders / number formats match
#    the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A8:Q8").Copy()
$ws.Rows(9).Insert()
$ws.Application.CutCopyMode = $false

# Row height for the freshly inserted row (matches the rest of the table)
$ws.Rows(9).RowHeight = 25.5

# Re-number the sequence column (A) for the rows that shifted down because of
# the insert - they keep their own data, they just need their running index
# bumped by one.
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# ---------------------------------------------------------------------------
# 2. Recreate the merged-cell layout for the new row 9 (same pattern as every
#    other data row: A:B, C:G, H:K, L:M, N:O).
# ---------------------------------------------------------------------------
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# ---------------------------------------------------------------------------
# 3. Fill in the values for the new item: EMEREST 4MG/2ML 5 AMPOULES
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = 3

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "EMEREST 4MG/2ML 5 AMPOULES"

$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "0:1"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"

$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "112.50"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "22.5000"

$ws.Range("Q9").NumberFormat = "@"
$ws.Range("Q9").Value = "0:1"

# Re-apply the formatting (font/fill/border/number format) from the sibling
# row so every cell in row 9 ends up using the exact same style as the rest
# of the table (setting NumberFormat above creates its own style ids; this
# normalizes them back).
$ws.Range("A10:Q10").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Update the totals row - the "sell price" column total grows by the new
#    item's sell price (267 + 22.50 = 289.50). After the insert this row is
#    now row 13.
# ---------------------------------------------------------------------------
$ws.Range("P13").Value = 289.5

# ---------------------------------------------------------------------------
# 5. Update the printed timestamp in the footer (now row 14).
# ---------------------------------------------------------------------------
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Monday, 11 August, 2025 10:35 AM"

Write-Host "Edit applied"
